# Generate Report for Handback
# Replace the two handback file identifiers (and their dependent timestamp /
# xliff-name columns) across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newId1 = "330fd509-94e5-4011-aec3-2893a3af9f11"
$newId2 = "ffff5b3b7d53-d086-4552-ab99-5fec2fd5f4ff"
$newHash = "edf0d6f0e5e0ca4dce42faf24f1506a81a7a3db3"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newId1.md"
$ws.Range("B2").Value = "e2e\$newId1.md"
$ws.Range("G2").Value = "2016-08-18 05:02:40"

$ws.Range("A3").Value = "$newId2.md"
$ws.Range("B3").Value = "e2e\$newId2.md"
$ws.Range("G3").Value = "2016-08-18 05:02:40"

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 2) {
        $hl.TextToDisplay = "e2e\$newId1.md"
    }
    elseif ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 2) {
        $hl.TextToDisplay = "e2e\$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newId1.md"
$ws.Range("G2").Value = "$newId1.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-18 05:02:35"
$ws.Range("I2").Value = "$newId1.md"
$ws.Range("J2").Value = "$newId1.$newHash.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-18 05:02:52"

$ws.Range("A3").Value = "$newId2.md"
$ws.Range("G3").Value = "$newId1.$newHash.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-18 05:02:35"
$ws.Range("I3").Value = "$newId2.md"
$ws.Range("J3").Value = "$newId1.$newHash.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-18 05:02:52"

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "$newId1.md"
    }
    elseif ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 9) {
        $hl.TextToDisplay = "$newId1.md"
    }
    elseif ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "$newId2.md"
    }
    elseif ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 9) {
        $hl.TextToDisplay = "$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newId1.md"
$ws.Range("G2").Value = "$newId1.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-18 05:02:40"
$ws.Range("I2").Value = "$newId1.md"
$ws.Range("J2").Value = "$newId1.$newHash.de-de.xlf"
$ws.Range("K2").Value = "2016-08-18 05:03:04"

$ws.Range("A3").Value = "$newId2.md"
$ws.Range("G3").Value = "$newId1.$newHash.de-de.xlf"
$ws.Range("H3").Value = "2016-08-18 05:02:40"
$ws.Range("I3").Value = "$newId2.md"
$ws.Range("J3").Value = "$newId1.$newHash.de-de.xlf"
$ws.Range("K3").Value = "2016-08-18 05:03:04"

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "$newId1.md"
    }
    elseif ($hl.Range.Row -eq 2 -and $hl.Range.Column -eq 9) {
        $hl.TextToDisplay = "$newId1.md"
    }
    elseif ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "$newId2.md"
    }
    elseif ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 9) {
        $hl.TextToDisplay = "$newId2.md"
    }
}
